$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New merged / centered header row (row 7) ---------------------------
$ws.Range("A7:C7").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A7").Value = "nộp bài muộn  (BTVN: 11/08/2020)"
$ws.Range("A7:C7").Merge()

# --- New name rows (8-17), column A formatted like the existing names ---
$ws.Range("A8:A17").Value = "placeholder"
$ws.Range("A1").Copy()
$ws.Range("A8:A17").PasteSpecial(-4122)          # xlPasteFormats

$ws.Range("A8").Value = "Ngô Tiến Thuận"
$ws.Range("A9").Value = "Đoàn Mạnh Cường"
$ws.Range("A10").Value = "Minh Bếu"
$ws.Range("A11").Value = "Nguyễn Trung Hiếu"
$ws.Range("A12").Value = "Hùng Phùng"
$ws.Range("A13").Value = "Vũ Thành Long"
$ws.Range("B13").Value = "10k"
$ws.Range("A14").Value = "Hiếu Bùi "
$ws.Range("A15").Value = "Hiếu Nguyễn"
$ws.Range("A16").Value = "Nguyen Tien Dung"
$ws.Range("A17").Value = "Lê Hiếu"

# --- Restore the reported selection --------------------------------------
$ws.Range("F12").Select()
